# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the per-language report sheets, adds hyperlinks to the newly
# populated target-file cells, and flips the overall status text from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdTarget49 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95fd08a224e94b37e816f4950301a3a8a0959bc4/e2e/49f9aacf-a83e-4802-83e9-ab74ce64ea86.md"
$mdTarget8f = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95fd08a224e94b37e816f4950301a3a8a0959bc4/e2e/8f756129-57d4-460f-8243-62527af05f51.md"

$statusText = "Handed back: in sync with en-US"

function Update-LangSheet($sheetName, $xlfDate, $xlf49, $xlf8f) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column now reflects the handback.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I) now points at the handed-back markdown file.
    $ws.Range("I2").Value = "49f9aacf-a83e-4802-83e9-ab74ce64ea86.md"
    $ws.Range("I3").Value = "8f756129-57d4-460f-8243-62527af05f51.md"

    $ws.Hyperlinks.Add($ws.Range("I2"), $mdTarget49, "", "", "49f9aacf-a83e-4802-83e9-ab74ce64ea86.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdTarget8f, "", "", "8f756129-57d4-460f-8243-62527af05f51.md") | Out-Null

    # Give those two cells the same hyperlink look as column A (applied
    # after Hyperlinks.Add so it wins over the default theme hyperlink font).
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = $true
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File (J) / Latest Handback DateTime (K).
    $ws.Range("J2").Value = $xlf49
    $ws.Range("J3").Value = $xlf8f
    $ws.Range("K2").Value = $xlfDate
    $ws.Range("K3").Value = $xlfDate

    # Widen the columns that now hold longer text, mirroring AutoFit.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Update-LangSheet "zh-cn" "2016-08-22 11:53:15" `
    "49f9aacf-a83e-4802-83e9-ab74ce64ea86.a2711cb97cecdcb1311cfd508a2c39b9fd243447.zh-cn.xlf" `
    "8f756129-57d4-460f-8243-62527af05f51.eac5e021f9e34bb47d15ce40bcf58813bd268230.zh-cn.xlf"

Update-LangSheet "de-de" "2016-08-22 11:53:22" `
    "49f9aacf-a83e-4802-83e9-ab74ce64ea86.a2711cb97cecdcb1311cfd508a2c39b9fd243447.de-de.xlf" `
    "8f756129-57d4-460f-8243-62527af05f51.eac5e021f9e34bb47d15ce40bcf58813bd268230.de-de.xlf"

# The Overview sheet shares the same "status" shared string for both
# languages; refresh it explicitly too (it also widens columns E/F).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527
